# Update countries & provincias Spain
# - Ucrania overtakes Ecuador in "Casos totales" ranking, so the two rows
#   swap country labels (new Ucrania figures move into the former Ecuador
#   row, and the old Ecuador figures slide down into the former Ucrania row).
# - Refresh several other countries' daily COVID figures.
# - Bump the "Datos actualizados" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 08:12"

# --- Row 19: Pakistan ---------------------------------------------------
$ws.Cells.Item(19,2).Value = 294638
$ws.Cells.Item(19,3).Value = 445
$ws.Cells.Item(19,4).Value = 279561
$ws.Cells.Item(19,5).Value = 8803
$ws.Cells.Item(19,7).Value = 7
$ws.Cells.Item(19,8).Value = 6274

# --- Row 30: now Ucrania (overtook Ecuador) -----------------------------
$ws.Cells.Item(30,1).Value = "Ucrania"
$ws.Cells.Item(30,2).Value = 112059
$ws.Cells.Item(30,3).Value = 1974
$ws.Cells.Item(30,4).Value = 54217
$ws.Cells.Item(30,5).Value = 55439
$ws.Cells.Item(30,7).Value = 49
$ws.Cells.Item(30,8).Value = 2403

# --- Row 31: now Ecuador (carries the old row-30 figures) ---------------
$ws.Cells.Item(31,1).Value = "Ecuador"
$ws.Cells.Item(31,2).Value = 110549
$ws.Cells.Item(31,3).Value = 0
$ws.Cells.Item(31,4).Value = 95097
$ws.Cells.Item(31,5).Value = 9042
$ws.Cells.Item(31,7).Value = 0
$ws.Cells.Item(31,8).Value = 6410

# --- Row 32: Israel -------------------------------------------------------
$ws.Cells.Item(32,2).Value = 108964
$ws.Cells.Item(32,3).Value = 561
$ws.Cells.Item(32,4).Value = 87011
$ws.Cells.Item(32,5).Value = 21078

# --- Row 57: Kirguistan ----------------------------------------------------
$ws.Cells.Item(57,2).Value = 43459
$ws.Cells.Item(57,3).Value = 101
$ws.Cells.Item(57,4).Value = 37492
$ws.Cells.Item(57,5).Value = 4910

# --- Row 62: Uzbekistan -----------------------------------------------------
$ws.Cells.Item(62,2).Value = 40195
$ws.Cells.Item(62,3).Value = 231
$ws.Cells.Item(62,5).Value = 3338
$ws.Cells.Item(62,7).Value = 4
$ws.Cells.Item(62,8).Value = 295

# --- Row 73: El Salvador ------------------------------------------------------
$ws.Cells.Item(73,2).Value = 25284
$ws.Cells.Item(73,3).Value = 144
$ws.Cells.Item(73,4).Value = 13291
$ws.Cells.Item(73,5).Value = 11299
$ws.Cells.Item(73,7).Value = 7
$ws.Cells.Item(73,8).Value = 694

# --- Row 180: Mauricio -----------------------------------------------------
$ws.Cells.Item(180,2).Value = 354
$ws.Cells.Item(180,5).Value = 9

# --- Row 188: Butan -------------------------------------------------------
$ws.Cells.Item(188,2).Value = 183
$ws.Cells.Item(188,3).Value = 10
$ws.Cells.Item(188,5).Value = 65

# --- Row 201: Fiyi ---------------------------------------------------------
$ws.Cells.Item(201,5).Value = 3
$ws.Cells.Item(201,7).Value = 1
$ws.Cells.Item(201,8).Value = 2
